# Generate Report for Handback
# Update the generated timestamps on the handback-status report:
#  - Overview!G2            (Latest HO Xliff Generate Date)
#  - zh-cn!H2                (Correspond Handoff Datetime)
#  - zh-cn!K2                (Correspond Handback DateTime)
#  - de-de!K2                (Correspond Handback DateTime)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-18 12:50:22"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-18 12:50:11"
$wsZhCn.Range("K2").Value = "2016-10-18 12:50:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-10-18 12:51:14"
